$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Vcam1"
$ws.Cells.Item(2,3).Value = "Itgad"
$ws.Cells.Item(2,4).Value = "M1"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 14.972416
$ws.Cells.Item(2,8).Value = 29.944832
$ws.Cells.Item(2,9).Value = 0.1033656722518705
$ws.Cells.Item(2,10).Value = 0.08332290573803899
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.1604463333333333
$ws.Cells.Item(2,14).Value = 0.481339
$ws.Cells.Item(2,15).Value = 0.4091752212750475
$ws.Cells.Item(2,16).Value = 0.4091752212750475
$ws.Cells.Item(2,17).Value = 2.402269248341333
$ws.Cells.Item(2,18).Value = 14.413615490048
$ws.Cells.Item(2,19).Value = 0.04229467181590314
$ws.Cells.Item(2,20).Value = 0.03409366839264203
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Vcam1"
$ws.Cells.Item(3,3).Value = "Itgad"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 14.972416
$ws.Cells.Item(3,8).Value = 29.944832
$ws.Cells.Item(3,9).Value = 0.1033656722518705
$ws.Cells.Item(3,10).Value = 0.08332290573803899
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.231675
$ws.Cells.Item(3,14).Value = 0.695025
$ws.Cells.Item(3,15).Value = 0.5908247787249524
$ws.Cells.Item(3,16).Value = 0.5908247787249525
$ws.Cells.Item(3,17).Value = 3.4687344768
$ws.Cells.Item(3,18).Value = 20.8124068608
$ws.Cells.Item(3,19).Value = 0.06107100043596733
$ws.Cells.Item(3,20).Value = 0.04922923734539696
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Vcam1"
$ws.Cells.Item(4,3).Value = "Itgad"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 40.41312266666666
$ws.Cells.Item(4,8).Value = 121.239368
$ws.Cells.Item(4,9).Value = 0.2790017050179012
$ws.Cells.Item(4,10).Value = 0.3373542530344942
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.1604463333333333
$ws.Cells.Item(4,14).Value = 0.481339
$ws.Cells.Item(4,15).Value = 0.4091752212750475
$ws.Cells.Item(4,16).Value = 0.4091752212750475
$ws.Cells.Item(4,17).Value = 6.484137350416887
$ws.Cells.Item(4,18).Value = 58.357236153752
$ws.Cells.Item(4,19).Value = 0.1141605843868152
$ws.Cells.Item(4,20).Value = 0.1380370011334675
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Vcam1"
$ws.Cells.Item(5,3).Value = "Itgad"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 40.41312266666666
$ws.Cells.Item(5,8).Value = 121.239368
$ws.Cells.Item(5,9).Value = 0.2790017050179012
$ws.Cells.Item(5,10).Value = 0.3373542530344942
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.231675
$ws.Cells.Item(5,14).Value = 0.695025
$ws.Cells.Item(5,15).Value = 0.5908247787249524
$ws.Cells.Item(5,16).Value = 0.5908247787249525
$ws.Cells.Item(5,17).Value = 9.362710193799998
$ws.Cells.Item(5,18).Value = 84.26439174419998
$ws.Cells.Item(5,19).Value = 0.1648411206310859
$ws.Cells.Item(5,20).Value = 0.1993172519010267
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Vcam1"
$ws.Cells.Item(6,3).Value = "Itgad"
$ws.Cells.Item(6,4).Value = "M1"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 5.007042333333334
$ws.Cells.Item(6,8).Value = 15.021127
$ws.Cells.Item(6,9).Value = 0.03456732011577652
$ws.Cells.Item(6,10).Value = 0.04179699352128983
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.1604463333333333
$ws.Cells.Item(6,14).Value = 0.481339
$ws.Cells.Item(6,15).Value = 0.4091752212750475
$ws.Cells.Item(6,16).Value = 0.4091752212750475
$ws.Cells.Item(6,17).Value = 0.8033615832281111
$ws.Cells.Item(6,18).Value = 7.230254249053
$ws.Cells.Item(6,19).Value = 0.01414409085725826
$ws.Cells.Item(6,20).Value = 0.01710229407270549
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Vcam1"
$ws.Cells.Item(7,3).Value = "Itgad"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 5.007042333333334
$ws.Cells.Item(7,8).Value = 15.021127
$ws.Cells.Item(7,9).Value = 0.03456732011577652
$ws.Cells.Item(7,10).Value = 0.04179699352128983
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.231675
$ws.Cells.Item(7,14).Value = 0.695025
$ws.Cells.Item(7,15).Value = 0.5908247787249524
$ws.Cells.Item(7,16).Value = 0.5908247787249525
$ws.Cells.Item(7,17).Value = 1.160006532575
$ws.Cells.Item(7,18).Value = 10.440058793175
$ws.Cells.Item(7,19).Value = 0.02042322925851826
$ws.Cells.Item(7,20).Value = 0.02469469944858434
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Vcam1"
$ws.Cells.Item(8,3).Value = "Itgad"
$ws.Cells.Item(8,4).Value = "M1"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 12.521722
$ws.Cells.Item(8,8).Value = 37.565166
$ws.Cells.Item(8,9).Value = 0.08644671723528362
$ws.Cells.Item(8,10).Value = 0.1045268440862112
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.1604463333333333
$ws.Cells.Item(8,14).Value = 0.481339
$ws.Cells.Item(8,15).Value = 0.4091752212750475
$ws.Cells.Item(8,16).Value = 0.4091752212750475
$ws.Cells.Item(8,17).Value = 2.009064381919333
$ws.Cells.Item(8,18).Value = 18.081579437274
$ws.Cells.Item(8,19).Value = 0.03537185465324864
$ws.Cells.Item(8,20).Value = 0.04276979455815785
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Vcam1"
$ws.Cells.Item(9,3).Value = "Itgad"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 12.521722
$ws.Cells.Item(9,8).Value = 37.565166
$ws.Cells.Item(9,9).Value = 0.08644671723528362
$ws.Cells.Item(9,10).Value = 0.1045268440862112
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.231675
$ws.Cells.Item(9,14).Value = 0.695025
$ws.Cells.Item(9,15).Value = 0.5908247787249524
$ws.Cells.Item(9,16).Value = 0.5908247787249525
$ws.Cells.Item(9,17).Value = 2.90096994435
$ws.Cells.Item(9,18).Value = 26.10872949915001
$ws.Cells.Item(9,19).Value = 0.05107486258203497
$ws.Cells.Item(9,20).Value = 0.06175704952805334
$ws.Cells.Item(10,1).Value = "Neutro"
$ws.Cells.Item(10,2).Value = "Vcam1"
$ws.Cells.Item(10,3).Value = "Itgad"
$ws.Cells.Item(10,4).Value = "M1"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 11.74303333333333
$ws.Cells.Item(10,8).Value = 35.2291
$ws.Cells.Item(10,9).Value = 0.08107085287879548
$ws.Cells.Item(10,10).Value = 0.09802663038937569
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.1604463333333333
$ws.Cells.Item(10,14).Value = 0.481339
$ws.Cells.Item(10,15).Value = 0.4091752212750475
$ws.Cells.Item(10,16).Value = 0.4091752212750475
$ws.Cells.Item(10,17).Value = 1.884126640544445
$ws.Cells.Item(10,18).Value = 16.9571397649
$ws.Cells.Item(10,19).Value = 0.03317218416563796
$ws.Cells.Item(10,20).Value = 0.04011006818042009
$ws.Cells.Item(11,1).Value = "Neutro"
$ws.Cells.Item(11,2).Value = "Vcam1"
$ws.Cells.Item(11,3).Value = "Itgad"
$ws.Cells.Item(11,4).Value = "M2"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 11.74303333333333
$ws.Cells.Item(11,8).Value = 35.2291
$ws.Cells.Item(11,9).Value = 0.08107085287879548
$ws.Cells.Item(11,10).Value = 0.09802663038937569
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.231675
$ws.Cells.Item(11,14).Value = 0.695025
$ws.Cells.Item(11,15).Value = 0.5908247787249524
$ws.Cells.Item(11,16).Value = 0.5908247787249525
$ws.Cells.Item(11,17).Value = 2.7205672475
$ws.Cells.Item(11,18).Value = 24.4851052275
$ws.Cells.Item(11,19).Value = 0.04789866871315751
$ws.Cells.Item(11,20).Value = 0.0579165622089556
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Vcam1"
$ws.Cells.Item(12,3).Value = "Itgad"
$ws.Cells.Item(12,4).Value = "M1"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 60.1916805
$ws.Cells.Item(12,8).Value = 120.383361
$ws.Cells.Item(12,9).Value = 0.4155477325003729
$ws.Cells.Item(12,10).Value = 0.3349723732305901
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.3333333333333333
$ws.Cells.Item(12,13).Value = 0.1604463333333333
$ws.Cells.Item(12,14).Value = 0.481339
$ws.Cells.Item(12,15).Value = 0.4091752212750475
$ws.Cells.Item(12,16).Value = 0.4091752212750475
$ws.Cells.Item(12,17).Value = 9.657534433396499
$ws.Cells.Item(12,18).Value = 57.945206600379
$ws.Cells.Item(12,19).Value = 0.1700318353961843
$ws.Cells.Item(12,20).Value = 0.1370623949376545
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Vcam1"
$ws.Cells.Item(13,3).Value = "Itgad"
$ws.Cells.Item(13,4).Value = "M2"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 60.1916805
$ws.Cells.Item(13,8).Value = 120.383361
$ws.Cells.Item(13,9).Value = 0.4155477325003729
$ws.Cells.Item(13,10).Value = 0.3349723732305901
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.231675
$ws.Cells.Item(13,14).Value = 0.695025
$ws.Cells.Item(13,15).Value = 0.5908247787249524
$ws.Cells.Item(13,16).Value = 0.5908247787249525
$ws.Cells.Item(13,17).Value = 13.9449075798375
$ws.Cells.Item(13,18).Value = 83.669445479025
$ws.Cells.Item(13,19).Value = 0.2455158971041885
$ws.Cells.Item(13,20).Value = 0.1979099782929356
